# The commit ("unify the conception of DataNode, DataTable, Entity") renames
# the worksheet that used to be called "Property1" to "DataNode". The rest of
# the underlying XML diff (new xr/xr2/xr3/xr16 revision namespaces, the
# fileVersion/rupBuild bump, the absPath pointing at a Windows path, the
# Calibri -> SimSun default-font swap and the row-height/dyDescent values it
# drags along, the "Normal" -> "常规" locale label, ...) is exactly the noise
# Excel stamps into a workbook when it is re-opened and re-saved on a
# different machine/locale/build - it is not a deliberate content edit, so it
# is not reproduced here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: Property1 -> DataNode
$ws.Name = "DataNode"

# The author's cursor ended up on F25 (bottom/frozen pane) when the file was
# saved; mirror that selection.
$ws.Range("F25").Select() | Out-Null
